$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 510
$ws.Range("F4").Value = 469
$ws.Range("F5").Value = 8824
$ws.Range("F7").Value = 11312
$ws.Range("F11").Value = 10
$ws.Range("F13").Value = 126
$ws.Range("F18").Value = 92
$ws.Range("F20").Value = 424
$ws.Range("F21").Value = 1916
$ws.Range("F22").Value = 728
$ws.Range("F23").Value = 654
$ws.Range("F24").Value = 364
$ws.Range("F25").Value = 297
$ws.Range("F27").Value = 613
$ws.Range("F29").Value = 1343
$ws.Range("F33").Value = 43
$ws.Range("F35").Value = 471
$ws.Range("F36").Value = 237
$ws.Range("F37").Value = 15
$ws.Range("F38").Value = 358
$ws.Range("F39").Value = 330
$ws.Range("F40").Value = 36
$ws.Range("F41").Value = 148
$ws.Range("F42").Value = 540
$ws.Range("F43").Value = 398
$ws.Range("F44").Value = 124
$ws.Range("F45").Value = 819
$ws.Range("F48").Value = 172
$ws.Range("F49").Value = 160

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F8").Value = 60
$ws.Range("F14").Value = 34
$ws.Range("F17").Value = 67
$ws.Range("F18").Value = 70
$ws.Range("F19").Value = 112
$ws.Range("F24").Value = 79
$ws.Range("F25").Value = 398

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 218
$ws.Range("F3").Value = 2841
$ws.Range("F5").Value = 217

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 218
$ws.Range("F4").Value = 510
$ws.Range("F6").Value = 217
$ws.Range("F7").Value = 8824
$ws.Range("F9").Value = 11313
$ws.Range("F12").Value = 10
$ws.Range("F14").Value = 126
$ws.Range("F18").Value = 424
$ws.Range("F19").Value = 1916
$ws.Range("F20").Value = 728
$ws.Range("F21").Value = 654
$ws.Range("F22").Value = 364
$ws.Range("F23").Value = 297
$ws.Range("F26").Value = 613
$ws.Range("F29").Value = 1343
$ws.Range("F31").Value = 34
$ws.Range("F34").Value = 67
$ws.Range("F36").Value = 238
$ws.Range("F37").Value = 358
$ws.Range("F39").Value = 148
$ws.Range("F40").Value = 540
$ws.Range("F41").Value = 398
$ws.Range("F42").Value = 124
$ws.Range("F45").Value = 398
$ws.Range("F48").Value = 172
$ws.Range("F49").Value = 160
